# Rotated several components to fit manufacturer's placement
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E38/E39: flip rotation sign (90 -> -90)
$ws.Range("E38").Value = -90
$ws.Range("E39").Value = -90

# E58..E60: re-rotate values
$ws.Range("E58").Value = 180
$ws.Range("E59").Value = 90
$ws.Range("E60").Value = 180

# Update the active selection to match the final cursor position left by the edit
$ws.Range("E61").Select()
